$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 0.092262
$ws.Range("H2").Value = 0.276786
$ws.Range("I2").Value = 0.674176843971804
$ws.Range("J2").Value = 0.6741768439718039
$ws.Range("M2").Value = 3.187097666666667
$ws.Range("N2").Value = 9.561293000000001
$ws.Range("O2").Value = 0.07264827682158614
$ws.Range("P2").Value = 0.09004258080252311
$ws.Range("Q2").Value = 0.294048004922
$ws.Range("R2").Value = 2.646432044298
$ws.Range("S2").Value = 0.04897778598756691
$ws.Range("T2").Value = 0.06070462294852117
$ws.Range("G3").Value = 0.092262
$ws.Range("H3").Value = 0.276786
$ws.Range("I3").Value = 0.674176843971804
$ws.Range("J3").Value = 0.6741768439718039
$ws.Range("O3").Value = 0.2149294457051256
$ws.Range("P3").Value = 0.2663903788010418
$ws.Range("Q3").Value = 0.8699390745879999
$ws.Range("R3").Value = 7.829451671291999
$ws.Range("S3").Value = 0.1449004553820908
$ws.Range("T3").Value = 0.1795942248445397
$ws.Range("G4").Value = 0.092262
$ws.Range("H4").Value = 0.276786
$ws.Range("I4").Value = 0.674176843971804
$ws.Range("J4").Value = 0.6741768439718039
$ws.Range("M4").Value = 2.693222
$ws.Range("N4").Value = 8.079666
$ws.Range("O4").Value = 0.06139063118282826
$ws.Range("P4").Value = 0.07608949737890039
$ws.Range("Q4").Value = 0.248482048164
$ws.Range("R4").Value = 2.236338433475999
$ws.Range("S4").Value = 0.04138814198027618
$ws.Range("T4").Value = 0.05129777720230791
$ws.Range("G5").Value = 0.092262
$ws.Range("H5").Value = 0.276786
$ws.Range("I5").Value = 0.674176843971804
$ws.Range("J5").Value = 0.6741768439718039
$ws.Range("M5").Value = 25.4243835
$ws.Range("N5").Value = 50.848767
$ws.Range("O5").Value = 0.5795359426364719
$ws.Range("P5").Value = 0.4788634979919735
$ws.Range("Q5").Value = 2.345704470477
$ws.Range("R5").Value = 14.074226822862
$ws.Range("S5").Value = 0.3907097127748811
$ws.Range("T5").Value = 0.3228386817695269
$ws.Range("G6").Value = 0.092262
$ws.Range("H6").Value = 0.276786
$ws.Range("I6").Value = 0.674176843971804
$ws.Range("J6").Value = 0.6741768439718039
$ws.Range("M6").Value = 3.136534000000001
$ws.Range("N6").Value = 9.409602000000001
$ws.Range("O6").Value = 0.07149570365398808
$ws.Range("P6").Value = 0.08861404502556121
$ws.Range("Q6").Value = 0.289382899908
$ws.Range("R6").Value = 2.604446099172
$ws.Range("S6").Value = 0.04820074784698906
$ws.Range("T6").Value = 0.05974153720690818
$ws.Range("I7").Value = 0.3258231560281961
$ws.Range("J7").Value = 0.3258231560281961
$ws.Range("M7").Value = 3.187097666666667
$ws.Range("N7").Value = 9.561293000000001
$ws.Range("O7").Value = 0.07264827682158614
$ws.Range("P7").Value = 0.09004258080252311
$ws.Range("Q7").Value = 0.1421105602248889
$ws.Range("R7").Value = 1.278995042024
$ws.Range("S7").Value = 0.02367049083401924
$ws.Range("T7").Value = 0.02933795785400194
$ws.Range("I8").Value = 0.3258231560281961
$ws.Range("J8").Value = 0.3258231560281961
$ws.Range("O8").Value = 0.2149294457051256
$ws.Range("P8").Value = 0.2663903788010418
$ws.Range("S8").Value = 0.07002899032303482
$ws.Range("T8").Value = 0.08679615395650209
$ws.Range("I9").Value = 0.3258231560281961
$ws.Range("J9").Value = 0.3258231560281961
$ws.Range("M9").Value = 2.693222
$ws.Range("N9").Value = 8.079666
$ws.Range("O9").Value = 0.06139063118282826
$ws.Range("P9").Value = 0.07608949737890039
$ws.Range("Q9").Value = 0.1200889734986667
$ws.Range("R9").Value = 1.080800761488
$ws.Range("S9").Value = 0.02000248920255209
$ws.Range("T9").Value = 0.02479172017659248
$ws.Range("I10").Value = 0.3258231560281961
$ws.Range("J10").Value = 0.3258231560281961
$ws.Range("M10").Value = 25.4243835
$ws.Range("N10").Value = 50.848767
$ws.Range("O10").Value = 0.5795359426364719
$ws.Range("P10").Value = 0.4788634979919735
$ws.Range("Q10").Value = 1.133656310676
$ws.Range("R10").Value = 6.801937864055999
$ws.Range("S10").Value = 0.1888262298615909
$ws.Range("T10").Value = 0.1560248162224465
$ws.Range("I11").Value = 0.3258231560281961
$ws.Range("J11").Value = 0.3258231560281961
$ws.Range("M11").Value = 3.136534000000001
$ws.Range("N11").Value = 9.409602000000001
$ws.Range("O11").Value = 0.07149570365398808
$ws.Range("P11").Value = 0.08861404502556121
$ws.Range("Q11").Value = 0.1398559600373334
$ws.Range("R11").Value = 1.258703640336
$ws.Range("S11").Value = 0.02329495580699903
$ws.Range("T11").Value = 0.02887250781865302
